$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: append more detail to the F25 (Status) note ---
$f25 = [string]$ws.Range("F25").Value2
$f25 += ".`nBefore AI processing, introduce a lightweight deduplication mechanism (e.g., store post_id, scraped_at in metadata)."
$ws.Range("F25").Value2 = $f25
$ws.Rows.Item(25).RowHeight = 57.6

# --- Row 26: add trailing space to the Task cell (C26) ---
$ws.Range("C26").Value2 = "Begin building the offline AI model pipeline "

# --- Row 27: fill in the Status cell (F27), which was previously empty ---
$ws.Range("F27").Value2 = "DONE: Old pipeline tested and worked perfectly, and the new branch is created for phase two"
$ws.Range("F27").HorizontalAlignment = -4131
$ws.Range("F27").VerticalAlignment = -4160
$ws.Range("F27").WrapText = $false

# --- Row 28: new Task/Context, plus Status cell ---
$ws.Range("C28").Value2 = "Make the super comprehensive roadmap for the full project"
$ws.Range("D28").Value2 = "After making 3 roadmaps combine the needed parts: Combine 3 ROADMAPS and FILE STRUTURES and MindMaps"
$ws.Range("D28").HorizontalAlignment = -4131
$ws.Range("D28").VerticalAlignment = -4160
$ws.Range("D28").WrapText = $false
$ws.Range("F28").Value2 = "DONE"
$ws.Range("F28").HorizontalAlignment = -4131
$ws.Range("F28").VerticalAlignment = -4160
$ws.Range("F28").WrapText = $false

# --- New rows 29-32: 31/7/2025(Onsite), continuing the Car Tracking Project log ---
$ws.Range("A29").Value2 = "31/7/2025(Onsite)"
$ws.Range("B29").Value2 = "Car Tracking Project"
$ws.Range("C29").Value2 = "Check how the medical o1 data set looked like so that you can format your data set in the same way"

$ws.Range("A30").Value2 = "31/7/2025(Onsite)"
$ws.Range("B30").Value2 = "Car Tracking Project"
$ws.Range("C30").Value2 = "Improve the file structure and readme file to include all the needed phases."

$ws.Range("A31").Value2 = "31/7/2025(Onsite)"
$ws.Range("B31").Value2 = "Car Tracking Project"
$ws.Range("C31").Value2 = "Generate a fake meta data for the branches so that you can test with later on"

$ws.Range("A32").Value2 = "31/7/2025(Onsite)"
$ws.Range("B32").Value2 = "Car Tracking Project"
$ws.Range("C32").Value2 = "Know the inputs and outputs of each phase and what happens at each phase"

# Match the styling (left/top aligned, no wrap) that the other A/B/C columns in this
# report use for plain single-line cells.
foreach ($r in 29..32) {
    foreach ($col in @("A","B","C")) {
        $rng = $ws.Range("$col$r")
        $rng.HorizontalAlignment = -4131
        $rng.VerticalAlignment = -4160
        $rng.WrapText = $false
    }
}

# --- Sheet view: scroll/zoom/selection state as left after the edit session ---
$ws.Application.ActiveWindow.Zoom = 103
$ws.Application.ActiveWindow.ScrollRow = 18
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A32").Select()
